$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, shifting existing rows 16-49 down to 17-50.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44498
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 300000001
$ws.Cells.Item(16, 7).Value = "Rabanito"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 50
$ws.Cells.Item(16, 11).Value = 7000
$ws.Cells.Item(16, 12).Value = 7000
$ws.Cells.Item(16, 13).Value = 7000
$ws.Cells.Item(16, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(16, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(16, 16).Value = 583
$ws.Cells.Item(16, 17).Value = 12
$ws.Cells.Item(16, 18).Value = "Hortaliza"
